$d = $word.ActiveDocument

# --- Paragraph 1: update the ID placeholder text and drop the trailing space run ---
$p1 = $d.Paragraphs(1)
$textLen = $p1.Range.Text.Length
$spaceRange = $d.Range($textLen - 2, $textLen - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

$d.Content.Find.Execute("**ID__AFFARS_5311_topic_3__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5311_103__ID**", 2)

# --- Paragraph 1: add a 5-twip paragraph border on all sides and widen the left indent ---
$p1 = $d.Paragraphs(1)
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25
